$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -1.523834722534454
$ws.Range("B2").Value = -4.988536054339994

$ws.Range("A3").Value = -0.5262281452645926
$ws.Range("B3").Value = 0.8912399176695573

$ws.Range("A4").Value = 0.9808900932775177
$ws.Range("B4").Value = -3.508034782370383

$ws.Range("A5").Value = 0.7146042438579954
$ws.Range("B5").Value = 0.5151879869131841

$ws.Range("A6").Value = -0.8051940605643112
$ws.Range("B6").Value = -1.712685417610642

$ws.Range("A7").Value = -0.1214127429496217
$ws.Range("B7").Value = 0.9957366901830694
